$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 27 and 28: Toncoin/Cosmos swap places with updated values
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.69%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.03%  '

# Remaining cell value updates (price + volume columns)
$ws.Range('D2').Value = '43.303.90'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '2.348.71'
$ws.Range('E3').Value = '  +4.00%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.649'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '231.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.37'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.12%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.457'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0942'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.88'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '2.692.96'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.105'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.61%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '2.345.38'
$ws.Range('E18').Value = '  +3.95%  '
$ws.Range('D19').Value = '43.265.18'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('E24').Value = '  +20.27%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '175.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.45%  '
$ws.Range('E31').Value = '  +8.58%  '
$ws.Range('E32').Value = '  -7.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.126'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0691'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.94%  '
$ws.Range('E37').Value = '  +8.60%  '
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0253'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('E42').Value = '  +8.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.93%  '
$ws.Range('E44').Value = '  +7.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '98.65'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0946'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('D49').Value = '1.437.42'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '2.565.95'
$ws.Range('E50').Value = '  +3.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000203'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.88%  '
